{"js": "// Fix word order in the recurring French \"Dates \u00e0 utiliser...\" campaign\n// sentence: \"Campagne Constellation d'Hercule 2022\" ->\n// \"Campagne 2022 Constellation d'Hercule\". The document repeats this\n// exact sentence in four places; replace all of them.\nconst oldText =\n  \"Dates \u00e0 utiliser pour la Campagne Constellation d'Hercule 2022: 13-22 juin, 12-21 juillet, 10-19 ao\u00fbt\";\nconst newText =\n  \"Dates \u00e0 utiliser pour la Campagne 2022 Constellation d'Hercule: 13-22 juin, 12-21 juillet, 10-19 ao\u00fbt\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix word order in the recurring French \"Dates \u00e0 utiliser...\" campaign\n# sentence: \"Campagne Constellation d'Hercule 2022\" ->\n# \"Campagne 2022 Constellation d'Hercule\". The document repeats this\n# exact sentence in four places; replace all of them.\n#\n# NOTE: we deliberately set Range.Text directly (rather than driving the\n# replacement through Find.Execute's ReplaceWith/Replacement.Text) because\n# Find/Replace runs the text through Word's AutoCorrect \"smart quotes\"\n# pass, which would turn the straight apostrophe in \"d'Hercule\" into a\n# curly one. Assigning Range.Text performs a plain text substitution and\n# keeps the original straight apostrophe, matching the source document.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Campagne Constellation d'Hercule 2022\"\n$newText = \"Campagne 2022 Constellation d'Hercule\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $oldText\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0          # wdFindStop - do not wrap, so the loop terminates\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n\nwhile ($rng.Find.Execute()) {\n  $rng.Text = $newText\n  $rng.Collapse(0)          # wdCollapseEnd - continue searching after this match\n  $rng.End = $d.Content.End\n}\n"}
